$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.892.69"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "'2.489.28"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'311.63"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'95.50"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").Value = "'0.556"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").Value = "'34.18"
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "'7.04"
$ws.Range("E13").Value = "  -3.87%  "
$ws.Range("D14").Value = "'2.870.50"
$ws.Range("D15").Value = "'2.497.32"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "'14.82"
$ws.Range("E16").Value = "  -5.55%  "
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "'41.905.27"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  -5.02%  "
$ws.Range("D20").Value = "'0.0₃0923"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").Value = "'69.51"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'237.71"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "'2.80"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").Value = "'1.94"
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").Value = "'36.71"
$ws.Range("E30").Value = "  -6.15%  "
$ws.Range("D31").Value = "'154.92"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -8.36%  "
$ws.Range("D35").Value = "'0.0761"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("D36").Value = "'3.05"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("D37").Value = "'17.30"
$ws.Range("E37").Value = "  -4.12%  "
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").Value = "'4.03"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("D42").Value = "'21.36"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "'2.000.70"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").Value = "'0.0287"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("D46").Value = "'3.10"
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'2.728.59"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").Value = "'77.38"
$ws.Range("E49").Value = "  -4.08%  "
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("D51").Value = "'98.38"
$ws.Range("E51").Value = "  -2.93%  "
